$d = $word.ActiveDocument

# Target content: 5 paragraphs, each made of 5 single-digit runs.
$rows = @(
    @("1","5","6","1","3"),
    @("2","6","2","4","2"),
    @("4","5","4","5","6"),
    @("2","2","3","5","2"),
    @("4","2","1","3","1")
)

$body = ""
foreach ($row in $rows) {
    $runs = ""
    foreach ($digit in $row) {
        $runs += "<w:r><w:t>$digit</w:t></w:r>"
    }
    $body += "<w:p>$runs</w:p>"
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       "<w:body>$body</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# Target the first (only) paragraph, excluding its trailing paragraph mark,
# so InsertXML replaces the run content without leaving a stray empty paragraph.
$full = $d.Paragraphs(1).Range
$target = $d.Range($full.Start, $full.End - 1)
$target.InsertXML($xml)

Write-Output $d.Content.Text
